$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: invoice_id text "4" -> "6" (kept as text, like the original inline string)
$ws.Cells.Item(2,1).Value = "'6"
$ws.Cells.Item(2,1).Style = "Normal"

# G2/H2/I2: updated product price / unit / invoice total
$ws.Range("G2").Value = 1002.21
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1182.6078

# Row 3 (the "bike" line item) is removed entirely
$ws.Rows("3").Delete()
